$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-02 Monday", "2024-12-03 Tuesday"),
    @("59×48=2832", "27×57=1539"),
    @("15×28=420", "94×32=3008"),
    @("86×24=2064", "19×32=608"),
    @("52×26=1352", "26×30=780"),
    @("13×37=481", "37×43=1591"),
    @("89×65=5785", "88×30=2640"),
    @("40×34=1360", "73×70=5110"),
    @("29×66=1914", "43×66=2838"),
    @("41×14=574", "41×54=2214"),
    @("64×38=2432", "78×35=2730"),
    @("70×21=1470", "81×82=6642"),
    @("46×54=2484", "12×86=1032"),
    @("47×88=4136", "91×71=6461"),
    @("67×19=1273", "74×93=6882"),
    @("24×28=672", "48×51=2448"),
    @("74×39=2886", "92×93=8556"),
    @("94×39=3666", "18×54=972"),
    @("18×33=594", "86×41=3526"),
    @("60×79=4740", "16×13=208"),
    @("59×15=885", "20×30=600"),
    @("75×60=4500", "91×49=4459"),
    @("32×51=1632", "34×26=884"),
    @("49×41=2009", "99×77=7623"),
    @("87×37=3219", "74×21=1554"),
    @("69×47=3243", "70×18=1260")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
